$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Sheet "Estimated" (sheet1): selection only moves, no data changes
# ------------------------------------------------------------------
$wsEstimated = $wb.Worksheets.Item("Estimated")
$wsEstimated.Activate() | Out-Null
$wsEstimated.Range("E2").Select() | Out-Null

# ------------------------------------------------------------------
# Sheet "current_total_expense_base" (sheet2): data + selection change,
# loses tabSelected
# ------------------------------------------------------------------
$wsBase = $wb.Worksheets.Item("current_total_expense_base")

$wsBase.Range("E2").Value = 0
$wsBase.Range("F2").Value = 1

$wsBase.Range("B3").Value = 0
$wsBase.Range("D3").Value = 0
$wsBase.Range("E3").Value = 0

$wsBase.Range("B4").Value = 0
$wsBase.Range("D4").Value = 0
$wsBase.Range("E4").Value = 0

$wsBase.Range("B5").Value = 0
$wsBase.Range("D5").Value = 0
$wsBase.Range("E5").Value = 0
$wsBase.Range("F5").Value = 1

$wsBase.Range("E6").Value = 0
$wsBase.Range("F6").Value = 1

$wsBase.Range("E7").Value = 0
$wsBase.Range("F7").Value = 1

$wsBase.Range("E8").Value = 0
$wsBase.Range("F8").Value = 1

$wsBase.Range("B9").Value = 0
$wsBase.Range("D9").Value = 0
$wsBase.Range("E9").Value = 0
$wsBase.Range("F9").Value = 1

$wsBase.Range("B10").Value = 0
$wsBase.Range("D10").Value = 0
$wsBase.Range("E10").Value = 0
$wsBase.Range("F10").Value = 1

$wsBase.Range("E11").Value = 0

$wsBase.Range("B12").Value = 0
$wsBase.Range("D12").Value = 0
$wsBase.Range("E12").Value = 0
$wsBase.Range("F12").Value = 1

$wsBase.Range("B13").Value = 0
$wsBase.Range("D13").Value = 0
$wsBase.Range("E13").Value = 0
$wsBase.Range("F13").Value = 1

$wsBase.Range("B14").Value = 0
$wsBase.Range("D14").Value = 0
$wsBase.Range("E14").Value = 0
$wsBase.Range("F14").Value = 1

$wsBase.Range("E15").Value = 0
$wsBase.Range("F15").Value = 1

$wsBase.Range("B16").Value = 0
$wsBase.Range("D16").Value = 0
$wsBase.Range("E16").Value = 0
$wsBase.Range("F16").Value = 1

$wsBase.Range("B17").Value = 0
$wsBase.Range("D17").Value = 0
$wsBase.Range("E17").Value = 0
$wsBase.Range("F17").Value = 1

$wsBase.Range("D18").Value = 0
$wsBase.Range("E18").Value = 0
$wsBase.Range("F18").Value = 1

$wsBase.Activate() | Out-Null
$wsBase.Range("M12").Select() | Out-Null

# ------------------------------------------------------------------
# Sheet "planned_estimated_cost_v1" (sheet6): quantities filled in,
# formulas added for the cost + grand total columns
# ------------------------------------------------------------------
$wsPlanned = $wb.Worksheets.Item("planned_estimated_cost_v1")

$wsPlanned.Range("B2").Value = 1
$wsPlanned.Range("D2").Value = 3000
$wsPlanned.Range("E2").Formula = "=SUM(D2:D17)"

$wsPlanned.Range("B3").Value = 8
$wsPlanned.Range("B4").Value = 150
$wsPlanned.Range("B5").Value = 1
$wsPlanned.Range("B6").Value = 1
$wsPlanned.Range("B7").Value = 1
$wsPlanned.Range("B8").Value = 1
$wsPlanned.Range("B9").Value = 1
$wsPlanned.Range("B10").Value = 1
$wsPlanned.Range("B11").Value = 12
$wsPlanned.Range("B12").Value = 4
$wsPlanned.Range("B13").Value = 2
$wsPlanned.Range("B14").Value = 1
$wsPlanned.Range("B15").Value = 1
$wsPlanned.Range("B16").Value = 1
$wsPlanned.Range("B17").Value = 1

$wsPlanned.Range("D3:D16").Formula = "=C3*B3"
$wsPlanned.Range("D17").Value = 1000

$wsPlanned.Range("E3").Value = 23784
$wsPlanned.Range("E4").Value = 23784
$wsPlanned.Range("E5").Value = 23784
$wsPlanned.Range("E6").Value = 23784
$wsPlanned.Range("E7").Value = 23784
$wsPlanned.Range("E8").Value = 23784
$wsPlanned.Range("E9").Value = 23784
$wsPlanned.Range("E10").Value = 23784
$wsPlanned.Range("E11").Value = 23784
$wsPlanned.Range("E12").Value = 23784
$wsPlanned.Range("E13").Value = 23784
$wsPlanned.Range("E14").Value = 23784
$wsPlanned.Range("E15").Value = 23784
$wsPlanned.Range("E16").Value = 23784
$wsPlanned.Range("E17").Value = 23784

$wsPlanned.Activate() | Out-Null
$wsPlanned.Range("E3:E17").Select() | Out-Null

# ------------------------------------------------------------------
# Sheet "current_total_expense_v1" (sheet5): data changes, and it
# ends up the final active sheet / tabSelected (activeTab="4", 0-based)
# so this block runs last.
# ------------------------------------------------------------------
$wsV1 = $wb.Worksheets.Item("current_total_expense_v1")
$wsV1.Range("B2").Value = 21000
$wsV1.Range("C2").Value = 30
$wsV1.Range("D2").Value = 0

$wsV1.Activate() | Out-Null
$wsV1.Range("D3").Select() | Out-Null
